# Vise manje gotov player stats page, dodaj DATE FILTER
#
# Slide 3 ("Graph: team points trend") has a note textbox
# ("TekstniOkvir 7") that currently reads "DODAJ FILTER ZA STATISTIKU".
# Widen the textbox so the extra text fits, and extend the note text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$noteShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "TekstniOkvir 7") {
        $noteShape = $sh
        break
    }
}

# Widen the textbox (height stays the same, autofit keeps it a single line)
$noteShape.Width = 445.2

# Update the note text
$noteShape.TextFrame.TextRange.Text = "DODAJ FILTER ZA STATISTIKU, u smislu sta ide na y os"
